$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Phase"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Units"
$ws.Range("E1").Value = "Note"
# F1 "Distribution Type" is unchanged
$ws.Range("G1").Value = "Parameter 1"
$ws.Range("H1").Value = "Parameter 2"
$ws.Range("I1").Value = "Parameter 3"
$ws.Range("J1").Value = "Parameter 4"
$ws.Range("K1").Value = "Parameter 5"
$ws.Range("L1").Value = "Parameter 6"
$ws.Range("M1").Value = "Parameter 7"
$ws.Range("N1").Value = "Lower Limit"
$ws.Range("O1").Value = "Upper Limit"
$ws.Range("P1").Value = "Step"

# --- Data rows (2-7): Name/Description/Units shift left one column, ---
# --- new Phase column A populated, Units column (old E) cleared,    ---
# --- and new Lower Limit/Upper Limit/Step columns (N:P) populated.  ---
$names = @("Test Parameter 1", "Test Parameter 2", "Test Parameter 3", "Test Parameter 4", "Test Parameter 5", "Test Parameter 6")
$descriptions = @("The first test parameter", "The second test parameter", "The third test parameter", "The fourth test parameter", "The fifth test parameter", "The sixth test parameter")

for ($i = 0; $i -lt 6; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = "Indoor"
    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = $descriptions[$i]
    $ws.Cells.Item($row, 4).Value = "Some Units"
    $ws.Range("E" + $row).ClearContents()
    $ws.Cells.Item($row, 14).Value = 0
    $ws.Cells.Item($row, 15).Value = 100
    $ws.Cells.Item($row, 16).Value = 1
}

# --- Column widths: autofit the newly-populated parameter/limit columns ---
$ws.Range("G1:M7").EntireColumn.AutoFit()
$ws.Range("N1:O7").EntireColumn.AutoFit()
$ws.Range("P1:P7").EntireColumn.AutoFit()

# --- Selection matches the post-edit workbook state ---
$ws.Range("N6").Select()
